$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 258, shifting the existing data (rows 258-262) down to 259-263.
$ws.Rows.Item(258).Insert()

# Copy the date cell style (s="2", numFmt 165) from the row below onto the new row's D cell
$ws.Cells.Item(259, 4).Copy()
$ws.Cells.Item(258, 4).PasteSpecial(-4122)  # xlPasteFormats

# Populate the newly inserted row 258 with the new weekly record.
$ws.Cells.Item(258, 1).Value = 3
$ws.Cells.Item(258, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(258, 3).Value = "Coquimbo"
$ws.Cells.Item(258, 4).Value = 44890
$ws.Cells.Item(258, 5).Value = 5
$ws.Cells.Item(258, 6).Value = "Fruta"
$ws.Cells.Item(258, 7).Value = 100101
$ws.Cells.Item(258, 8).Value = "Berries"
$ws.Cells.Item(258, 9).Value = 100101001
$ws.Cells.Item(258, 10).Value = "Arándano (blue)"
$ws.Cells.Item(258, 11).Value = "Sin especificar"
$ws.Cells.Item(258, 12).Value = "Primera"
$ws.Cells.Item(258, 13).Value = 50
$ws.Cells.Item(258, 14).Value = 6000
$ws.Cells.Item(258, 15).Value = 6000
$ws.Cells.Item(258, 16).Value = 6000
$ws.Cells.Item(258, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(258, 18).Value = "Provincia de Linares"
$ws.Cells.Item(258, 19).Value = 3000
$ws.Cells.Item(258, 20).Value = 2
